$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.25
$ws.Range("E3").Value = 1.25
$ws.Range("E4").Value = 0.6
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 1.25

$ws.Range("J2:J6").ClearContents()

$ws.Range("E2").Select()
